$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.281.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "'1.681.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'218.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").Value = "'0.5276"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.93%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +1.56%  "
$ws.Range("D9").Value = "'0.06434"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("D10").Value = "'22.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.50%  "
$ws.Range("D11").Value = "'0.07497"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "'1.684.80"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "'0.5820"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").Value = "'0.000008494"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "'64.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").Value = "'26.322.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "'4.927"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("D21").Value = "'189.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "'6.210"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'144.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").Value = "'7.734"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").Value = "'0.1239"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.78%  "
$ws.Range("D27").Value = "'15.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("D28").Value = "'0.06642"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.16%  "
$ws.Range("D29").Value = "'1.360"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.02%  "
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("D32").Value = "'3.573"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("D33").Value = "'1.665"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("D34").Value = "'1.028"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("D35").Value = "'0.6208"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.99%  "
$ws.Range("D36").Value = "'2.397"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.01%  "
$ws.Range("D37").Value = "'2.700"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("E38").Value = "  +5.19%  "
$ws.Range("D39").Value = "'1.110.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.01%  "
$ws.Range("D40").Value = "'0.01625"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("D41").Value = "'0.8777"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.81%  "
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").Value = "'100.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("D44").Value = "'1.829.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("D45").Value = "'0.00000000109"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.92%  "
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("D47").Value = "'8.159"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("D48").Value = "'1.003"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("D49").Value = "'0.05268"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("D50").Value = "'0.4304"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").Value = "'6.056"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.95%  "
